# Edit the "KFF_colors" slide: append hex-code suffixes to a few rgb(...) labels,
# grow one label's textbox to fit its new (wrapped) text, and recolor a swatch
# rectangle to match its corrected hex value.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The color-swatch labels live inside the "Group 12" group shape (shape #2 on the slide).
$grp = $s.Shapes.Item(2)

# TextBox 1 (id=2): "rgb(10,63,100)" -> "rgb(10,63,100) #0A3F64"
$tb1 = $grp.GroupItems.Item(2)
$tb1.TextFrame.TextRange.Characters(4, 11).Text = "(10,63,100) #0A3F64"

# TextBox 8 (id=9): "rgb(12,91,152)" -> "rgb(12,91,152) #0C5B98"
$tb8 = $grp.GroupItems.Item(3)
$tb8.TextFrame.TextRange.Characters(4, 11).Text = "(12,91,152) #0C5B98"

# TextBox 9 (id=10): "rgb(19,120,189)" -> "rgb(19,120,189) #1378BD"
# and the textbox grows taller to accommodate the now-wrapped second line.
$tb9 = $grp.GroupItems.Item(4)
$tb9.TextFrame.TextRange.Characters(4, 12).Text = "(19,120,189) #1378BD"
$tb9.Height = 50.89221

# Rectangle 13 (id=14): recolor the swatch fill from 3E99D3 to the corrected 0C5B98.
$rect13 = $s.Shapes.Item(4)
$rect13.Fill.ForeColor.RGB = 9984780
